# The commit renames the embedded logo pictures that live in the
# document's headers/footers:
#   - the Pearson Edexcel logo (descr ...PearsonLogo.png) is renamed
#     from "image1.png" to "image2.png"
#   - the BTEC logo (descr "BTec_Logo-Orange") is renamed
#     from "image2.jpg" to "image1.jpg"
#
# These pictures live in wp:inline drawings inside the section
# headers/footers (not in the main body), so we have to walk
# Sections -> Headers/Footers -> Range.InlineShapes to reach them.

$d = $word.ActiveDocument

function Rename-LogoShape($shape) {
    if ($shape -eq $null) { return }

    $descr = $shape.AlternativeText
    $newName = $null

    if ($descr -eq "BTec_Logo-Orange") {
        $newName = "image1.jpg"
    } elseif ($descr -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $newName = "image2.png"
    }

    if ($newName -eq $null) { return }

    # Renaming an InlineShape reached directly through a
    # Header/Footer.Range can report a stale handle in some engines;
    # selecting the shape's range first and renaming it through the
    # Selection's InlineShapes collection is the robust path.
    $shape.Range.Select()
    $selShape = $word.Selection.InlineShapes.Item(1)
    $selShape.Name = $newName
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                Rename-LogoShape $shapes.Item($i)
            }
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                Rename-LogoShape $shapes.Item($i)
            }
        }
    }
}
